$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF"), matching the style of H1 ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-60: column I and J numeric values ---
$data = @(
    @(2,7,7),
    @(3,5,6),
    @(4,8,8),
    @(5,3,4),
    @(6,9,9),
    @(7,7,7),
    @(8,8,8),
    @(9,8,8),
    @(10,8,8),
    @(11,4,5),
    @(12,7,8),
    @(13,9,9),
    @(14,7,7),
    @(15,8,8),
    @(16,7,7),
    @(17,8,8),
    @(18,8,8),
    @(19,10,10),
    @(20,7,7),
    @(21,7,7),
    @(22,7,7),
    @(23,8,8),
    @(24,9,9),
    @(25,9,9),
    @(26,7,7),
    @(27,8,9),
    @(28,8,8),
    @(29,8,8),
    @(30,7,7),
    @(31,6,6),
    @(32,6,6),
    @(33,7,7),
    @(34,6,7),
    @(35,7,7),
    @(36,7,7),
    @(37,7,7),
    @(38,8,8),
    @(39,8,8),
    @(40,8,8),
    @(41,8,8),
    @(42,9,9),
    @(43,8,9),
    @(44,8,8),
    @(45,7,7),
    @(46,7,7),
    @(47,5,5),
    @(48,7,7),
    @(49,7,7),
    @(50,6,6),
    @(51,7,7),
    @(52,7,8),
    @(53,4,4),
    @(54,7,8),
    @(55,6,6),
    @(56,2,2),
    @(57,7,7),
    @(58,5,5),
    @(59,3,3),
    @(60,4,4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
